$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# BLEU score
$ws.Range("B11").Value = 0.1067931366085375

# Code BLEU
$ws.Range("B12").Value = 0.257156682201655
$ws.Range("C12").Value = "{'codebleu': 0.257156682201655, 'ngram_match_score': 0.10679313660853745, 'weighted_ngram_match_score': 0.1160108073879559, 'syntax_match_score': 0.54, 'dataflow_match_score': 0.26582278481012656}"

# Embeddings and Cosine similarity
$ws.Range("B13").Value = 0.933243252828546
